$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-17 from serial date 45224 (2023-10-25)
# to serial date 45233 (2023-11-03), matching the committed XML diff.
$ws.Range("C2:C17").Value = 45233
